$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update worker data in row 16 for the new account-statement entry ---
# (Old data: CC 1002059412 - ERIKA PATRICIA PEREZ MARTINEZ - period 2507)
# (New data: CC 1044926194 - ALVARO ENRIQUE MOSCOTE DE LA ROSA - period 2508)
$ws.Range("C16").Value = "1044926194"
$ws.Range("D16").Value = "ALVARO ENRIQUE MOSCOTE DE LA ROSA"
$ws.Range("E16").Value = "2508"

# --- Remove the extra period rows (17-20) that belonged to the old worker ---
# Only a single period (row 16) remains for the new worker/statement part 1
$ws.Range("17:20").EntireRow.Delete()

# --- Update the summary totals to reflect the single remaining period ---
$ws.Range("E11").Value = 56940
$ws.Range("F13").Value = 1

# --- Widen column D so the longer worker name fits (mirrors Excel's bestFit) ---
$ws.Columns("D:D").ColumnWidth = 37.3
